$d = $word.ActiveDocument

# Locate the "Ver no Jupiter ..." boilerplate paragraph with Find.
$searchRange = $d.Content
$found = $searchRange.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", $true, $false, $false,
                                    $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate anchor paragraph 'Ver no Jupiter...'"
}

# Map the matched range back to its 1-based paragraph index in the document.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $searchRange.Start -and $p.Range.End -ge $searchRange.End) {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -eq -1) {
    throw "Could not map found text back to a paragraph."
}

# The site footer consists of three consecutive paragraphs: a blank line,
# "Ver no Jupiter Salvar em pdf Salvar em docx", and the "(c) 2020 ..."
# copyright/credits line. All three are removed together, leaving the
# "LOQ4073: ..." paragraph directly followed by the page-break paragraph.
$startPara = $d.Paragraphs.Item($targetIndex - 1)
$endPara = $d.Paragraphs.Item($targetIndex + 1)

$delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$delRange.Delete()

Write-Output "Removed boilerplate paragraphs."
